# Auto-generated: apply cell-value updates per sheet, matching the target XML diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4133260
$ws.Range("I33").Value = 1190.4445
$ws.Range("K33").Value = 1190.4445
$ws.Range("M33").Value = -961.4445000000001
$ws.Range("H64").Value = 3997.5356
$ws.Range("I64").Value = 3903.7646
$ws.Range("J64").Value = 4142.4546
$ws.Range("K64").Value = 3903.7646
$ws.Range("L64").Value = 4142.4546
$ws.Range("M64").Value = -3655.7646
$ws.Range("N64").Value = -4638.4546
$ws.Range("H67").Value = 3997.5356
$ws.Range("I67").Value = 3903.7646
$ws.Range("J67").Value = 4142.4546
$ws.Range("K67").Value = 3903.7646
$ws.Range("L67").Value = 4142.4546
$ws.Range("M67").Value = -3045.7646
$ws.Range("N67").Value = -5858.4546
$ws.Range("H69").Value = 4466.6523
$ws.Range("I69").Value = 3800.25
$ws.Range("J69").Value = 5193.636
$ws.Range("K69").Value = 11400.75
$ws.Range("L69").Value = 15580.908
$ws.Range("M69").Value = -10526.75
$ws.Range("N69").Value = -17328.908
$ws.Range("H72").Value = 4466.6523
$ws.Range("I72").Value = 3800.25
$ws.Range("J72").Value = 5193.636
$ws.Range("K72").Value = 34202.25
$ws.Range("L72").Value = 46742.724
$ws.Range("M72").Value = -29834.25
$ws.Range("N72").Value = -55478.724
$ws.Range("H103").Value = 382982.4
$ws.Range("I103").Value = 630.2857
$ws.Range("J103").Value = 650628.9
$ws.Range("K103").Value = 1890.8571
$ws.Range("L103").Value = 1951886.7
$ws.Range("M103").Value = -1304.8571
$ws.Range("N103").Value = -1953058.7
$ws.Range("H116").Value = 7888.4707
$ws.Range("I116").Value = 16000.571
$ws.Range("J116").Value = 2210
$ws.Range("K116").Value = 16000.571
$ws.Range("L116").Value = 2210
$ws.Range("M116").Value = -12558.571
$ws.Range("N116").Value = -9094
$ws.Range("H118").Value = 657.2
$ws.Range("I118").Value = 521.5
$ws.Range("K118").Value = 1564.5
$ws.Range("M118").Value = 92.5
$ws.Range("H135").Value = 4034.5
$ws.Range("I135").Value = 4427.222
$ws.Range("J135").Value = 500
$ws.Range("K135").Value = 39844.998
$ws.Range("L135").Value = 4500
$ws.Range("M135").Value = -37309.998
$ws.Range("N135").Value = -9570
$ws.Range("H141").Value = 2085.625
$ws.Range("I141").Value = 1945
$ws.Range("J141").Value = 2695
$ws.Range("K141").Value = 5835
$ws.Range("L141").Value = 8085
$ws.Range("M141").Value = -655
$ws.Range("N141").Value = -18445

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4582.161
$ws.Range("I32").Value = 3320.4375
$ws.Range("K32").Value = 3320.4375
$ws.Range("M32").Value = -3033.4375
$ws.Range("H61").Value = 1402.7391
$ws.Range("I61").Value = 1330
$ws.Range("J61").Value = 1748.25
$ws.Range("K61").Value = 1330
$ws.Range("L61").Value = 1748.25
$ws.Range("M61").Value = -1118
$ws.Range("N61").Value = -2172.25
$ws.Range("H122").Value = 1604197.5
$ws.Range("I122").Value = 1833047.2
$ws.Range("K122").Value = 5499141.6
$ws.Range("M122").Value = -5496691.6
$ws.Range("H132").Value = 4709.1274
$ws.Range("I132").Value = 1355.3684
$ws.Range("J132").Value = 18869.445
$ws.Range("K132").Value = 4066.1052
$ws.Range("L132").Value = 56608.335
$ws.Range("M132").Value = -1536.1052
$ws.Range("N132").Value = -61668.335
$ws.Range("H136").Value = 1402.7391
$ws.Range("I136").Value = 1330
$ws.Range("J136").Value = 1748.25
$ws.Range("K136").Value = 3990
$ws.Range("L136").Value = 5244.75
$ws.Range("M136").Value = -1440
$ws.Range("N136").Value = -10344.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 33799.9
$ws.Range("J88").Value = 33799.9
$ws.Range("L88").Value = 33799.9
$ws.Range("N88").Value = -34611.9
$ws.Range("H91").Value = 33799.9
$ws.Range("J91").Value = 33799.9
$ws.Range("L91").Value = 33799.9
$ws.Range("N91").Value = -36607.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 415.36365
$ws.Range("I22").Value = 397.14285
$ws.Range("J22").Value = 447.25
$ws.Range("K22").Value = 397.14285
$ws.Range("L22").Value = 447.25
$ws.Range("M22").Value = -47.14285000000001
$ws.Range("N22").Value = -1147.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1635842.4
$ws.Range("I132").Value = 1955.8
$ws.Range("J132").Value = 1917547
$ws.Range("K132").Value = 17602.2
$ws.Range("L132").Value = 17257923
$ws.Range("M132").Value = -15072.2
$ws.Range("N132").Value = -17262983

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5249.4585
$ws.Range("I70").Value = 5167.243
$ws.Range("J70").Value = 5526
$ws.Range("K70").Value = 5167.243
$ws.Range("L70").Value = 5526
$ws.Range("M70").Value = -4897.243
$ws.Range("N70").Value = -6066
$ws.Range("H73").Value = 5249.4585
$ws.Range("I73").Value = 5167.243
$ws.Range("J73").Value = 5526
$ws.Range("K73").Value = 5167.243
$ws.Range("L73").Value = 5526
$ws.Range("M73").Value = -4231.243
$ws.Range("N73").Value = -7398
$ws.Range("H122").Value = 6755750.5
$ws.Range("I122").Value = 8103470
$ws.Range("J122").Value = 5557778
$ws.Range("K122").Value = 24310410
$ws.Range("L122").Value = 16673334
$ws.Range("M122").Value = -24307960
$ws.Range("N122").Value = -16678234

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2925827.2
$ws.Range("I22").Value = 12346097
$ws.Range("J22").Value = 2295.276
$ws.Range("K22").Value = 12346097
$ws.Range("L22").Value = 2295.276
$ws.Range("M22").Value = -12345802
$ws.Range("N22").Value = -2885.276
$ws.Range("H27").Value = 2925827.2
$ws.Range("I27").Value = 12346097
$ws.Range("J27").Value = 2295.276
$ws.Range("K27").Value = 12346097
$ws.Range("L27").Value = 2295.276
$ws.Range("M27").Value = -12345990
$ws.Range("N27").Value = -2509.276
$ws.Range("H93").Value = 25011330
$ws.Range("I93").Value = 17707.666
$ws.Range("J93").Value = 62501764
$ws.Range("K93").Value = 17707.666
$ws.Range("L93").Value = 62501764
$ws.Range("M93").Value = -16459.666
$ws.Range("N93").Value = -62504260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 22778.334
$ws.Range("J64").Value = 22778.334
$ws.Range("L64").Value = 22778.334
$ws.Range("N64").Value = -23274.334
$ws.Range("H67").Value = 22778.334
$ws.Range("J67").Value = 22778.334
$ws.Range("L67").Value = 22778.334
$ws.Range("N67").Value = -24494.334
$ws.Range("H136").Value = 3964.95
$ws.Range("I136").Value = 4949.9165
$ws.Range("J136").Value = 2487.5
$ws.Range("K136").Value = 14849.7495
$ws.Range("L136").Value = 7462.5
$ws.Range("M136").Value = -12299.7495
$ws.Range("N136").Value = -12562.5
